# UKTAG 20210528 MQTT.pptx - Sophie Germain slide bio edit
# - fixed a typo and refined the bio for Sophie Germain
#
# Target shape: "Rectangle 11" on slide 3 (the bullet-point bio textbox).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)

# ---------------------------------------------------------------------
# 1) Resize / reposition the textbox (it grows wider & moves up a bit so
#    the longer bio text still fits).
# ---------------------------------------------------------------------
$shp.Left   = 87.41779727559054
$shp.Top    = 126.26173028346457
$shp.Width  = 528.182098984252
$shp.Height = 269.0015718031496

$tr = $shp.TextFrame.TextRange

# ---------------------------------------------------------------------
# 2) Paragraph 3 ("Add to use a pseudonym to get in the math academic
#    circles (Antoine Auguste Le Blanc)") becomes "Corresponded profusely
#    with Carl Friedrich Gauss and others, using the pseudonym Antoine
#    Auguste Le Blanc" (with the name in italic blue).
# ---------------------------------------------------------------------

# "Add" -> "with"
$anchor = "Add"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Text = "with"

# Insert "Corresponded profusely " right before "with"
$anchor = "with"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.InsertBefore("Corresponded profusely ")

# " to use a " -> " Carl Friedrich Gauss and "
$anchor = " to use a "
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Text = " Carl Friedrich Gauss and "

# Insert "others, using the " right before "pseudonym"
$anchor = "pseudonym"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.InsertBefore("others, using the ")

# Remove " to get in the math academic circles (" that used to follow
# "pseudonym" (it becomes just a single space before "Antoine").
$anchor = " to get in the math academic circles ("
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Text = " "

# Drop the trailing ")" after "...Le Blanc"
$anchor = "Antoine Auguste Le Blanc)"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Text = "Antoine Auguste Le Blanc"

# Colour / italicize "Antoine Auguste Le Blanc"
$anchor = "Antoine Auguste Le Blanc"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Font.Italic = $true
$run.Font.Color.RGB = 12611584

# ---------------------------------------------------------------------
# 3) Paragraph 4 ("Probably saved Carl Friedrich Gauss from Napoleon's
#    troops…") becomes "...Pr. Gauss..." and gets a new clause appended.
# ---------------------------------------------------------------------

$anchor = " Carl Friedrich Gauss "
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1, $anchor.Length)
$run.Text = " Pr. Gauss "

$anchor = "troops"
$whole = $tr.Text
$pos = $whole.IndexOf($anchor)
$run = $tr.Characters($pos + 1 + $anchor.Length, 1)
$run.Text = "… and when he learned who she actually was."
